$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Fall2024/Spring2024/Summer2024 and Fall2025/Spring2025/Summer2025
# blocks entirely (rows 21 through 38).
$ws.Range("A21:F38").ClearContents()

# The courses previously listed under the Fall2023/Spring2023/Summer2023
# block (rows 13-14) move up into the Fall2022/Spring2022/Summer2022 block,
# and a couple of courses from the removed blocks join them too.

# Row 4: add Summer 2022 course (CPSC 4176) and update Spring 2022 course
# to CPSC 4148
$ws.Range("C4").Value = "CPSC 4148"
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "CPSC 4176"
$ws.Range("F4").Value = 3

# Row 5: Spring 2022 course becomes CPSC 4155
$ws.Range("C5").Value = "CPSC 4155"
$ws.Range("D5").Value = 3

# Row 6 (new): Fall 2022 course CPSC 3121, Spring 2022 course CPSC 4157
$ws.Range("A6").Value = "CPSC 3121"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "CPSC 4157"
$ws.Range("D6").Value = 3

# Row 7 (new): Fall 2022 course CPSC 3165, Spring 2022 course CPSC 4175
$ws.Range("A7").Value = "CPSC 3165"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "CPSC 4175"
$ws.Range("D7").Value = 3

# Row 8 (new): Fall 2022 course CPSC 4000, 0 credits
$ws.Range("A8").Value = "CPSC 4000"
$ws.Range("B8").Value = 0

# Row 9 (new): Fall 2022 course CPSC 4135
$ws.Range("A9").Value = "CPSC 4135"
$ws.Range("B9").Value = 3

# Clear out the old course rows (13-14) under the Fall2023 block since
# those courses have been moved up into the Fall2022 block above.
$ws.Range("A13:F14").ClearContents()
